$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Old and new JOIN blocks shared by every SQL query string on the sheet.
$oldBlock = "LEFT JOIN `n    df_participant prt ON std.id = prt.`"study.id`"`nLEFT JOIN `n    df_diagnoses dgn ON prt.id = dgn.`"participant.id`"`nLEFT JOIN `n    df_treatments trt ON prt.id = trt.`"participant.id`"`nLEFT JOIN `n    df_treatment_resp trr ON prt.id = trr.`"participant.id`"`nLEFT JOIN `n    df_survival srv ON prt.id = srv.`"participant.id`"`nLEFT JOIN `n    df_reference_files rfs ON std.id = rfs.`"study.id`""

$newBlock = "LEFT JOIN `n    df_participant prt ON std.study_id = prt.`"study.study_id`"`nLEFT JOIN `n    df_diagnoses dgn ON prt.participant_id = dgn.`"participant.participant_id`"`nLEFT JOIN `n    df_treatments trt ON prt.participant_id = trt.`"participant.participant_id`"`nLEFT JOIN `n    df_treatment_resp trr ON prt.participant_id = trr.`"participant.participant_id`"`nLEFT JOIN `n    df_survival srv ON prt.participant_id = srv.`"participant.participant_id`"`nLEFT JOIN `n    df_reference_files rfs ON std.study_id = rfs.`"study.study_id`""

function Update-QueryCell($cell) {
    $text = $cell.Value()
    if ($text -eq $null) { return }
    $text = $text.ToString()
    if ($text.Contains($oldBlock)) {
        $cell.Value = $text.Replace($oldBlock, $newBlock)
    }
}

# Column B (TabQuery) rows 2-7 and column C (StatQuery) row 2 hold the SQL text.
for ($r = 2; $r -le 7; $r++) {
    Update-QueryCell $ws.Cells.Item($r, 2)
}
Update-QueryCell $ws.Cells.Item(2, 3)

# Column C width change: from bestFit 60.83203125 to fixed 69.33203125 (no bestFit).
$ws.Columns.Item(3).ColumnWidth = 68.41666666666667
